{"js": "// Replace the division-problem answers in the table with the new set of\n// worked examples, per the commit \"Update master to output generated at c986bee\".\n// Each entry is unique in the document, so a direct text search + replace\n// for each pair is safe and unambiguous.\nconst replacements = [\n  [\"348\u00f76=58, 0\", \"724\u00f76=120, 4\"],\n  [\"409\u00f76=68, 1\", \"116\u00f79=12, 8\"],\n  [\"985\u00f79=109, 4\", \"754\u00f77=107, 5\"],\n  [\"918\u00f77=131, 1\", \"741\u00f72=370, 1\"],\n  [\"867\u00f74=216, 3\", \"891\u00f76=148, 3\"],\n  [\"345\u00f76=57, 3\", \"834\u00f78=104, 2\"],\n  [\"946\u00f75=189, 1\", \"769\u00f73=256, 1\"],\n  [\"944\u00f74=236, 0\", \"748\u00f77=106, 6\"],\n  [\"411\u00f72=205, 1\", \"299\u00f72=149, 1\"],\n  [\"546\u00f77=78, 0\", \"928\u00f79=103, 1\"],\n  [\"154\u00f74=38, 2\", \"461\u00f76=76, 5\"],\n  [\"807\u00f72=403, 1\", \"336\u00f74=84, 0\"],\n  [\"535\u00f77=76, 3\", \"929\u00f78=116, 1\"],\n  [\"959\u00f74=239, 3\", \"244\u00f75=48, 4\"],\n  [\"261\u00f74=65, 1\", \"973\u00f77=139, 0\"],\n  [\"653\u00f75=130, 3\", \"775\u00f74=193, 3\"],\n  [\"470\u00f74=117, 2\", \"719\u00f73=239, 2\"],\n  [\"486\u00f73=162, 0\", \"159\u00f76=26, 3\"],\n  [\"330\u00f75=66, 0\", \"309\u00f79=34, 3\"],\n  [\"319\u00f74=79, 3\", \"843\u00f78=105, 3\"],\n  [\"532\u00f79=59, 1\", \"361\u00f74=90, 1\"],\n  [\"269\u00f79=29, 8\", \"953\u00f78=119, 1\"],\n  [\"444\u00f74=111, 0\", \"763\u00f79=84, 7\"],\n  [\"389\u00f75=77, 4\", \"135\u00f72=67, 1\"],\n  [\"899\u00f77=128, 3\", \"549\u00f72=274, 1\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem answers in the table with the new set of\n# worked examples, per the commit \"Update master to output generated at c986bee\".\n# Each \"old\" value occurs exactly once in the document, so a direct\n# Find/Replace for each pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"348\u00f76=58, 0\", \"724\u00f76=120, 4\"),\n    @(\"409\u00f76=68, 1\", \"116\u00f79=12, 8\"),\n    @(\"985\u00f79=109, 4\", \"754\u00f77=107, 5\"),\n    @(\"918\u00f77=131, 1\", \"741\u00f72=370, 1\"),\n    @(\"867\u00f74=216, 3\", \"891\u00f76=148, 3\"),\n    @(\"345\u00f76=57, 3\", \"834\u00f78=104, 2\"),\n    @(\"946\u00f75=189, 1\", \"769\u00f73=256, 1\"),\n    @(\"944\u00f74=236, 0\", \"748\u00f77=106, 6\"),\n    @(\"411\u00f72=205, 1\", \"299\u00f72=149, 1\"),\n    @(\"546\u00f77=78, 0\", \"928\u00f79=103, 1\"),\n    @(\"154\u00f74=38, 2\", \"461\u00f76=76, 5\"),\n    @(\"807\u00f72=403, 1\", \"336\u00f74=84, 0\"),\n    @(\"535\u00f77=76, 3\", \"929\u00f78=116, 1\"),\n    @(\"959\u00f74=239, 3\", \"244\u00f75=48, 4\"),\n    @(\"261\u00f74=65, 1\", \"973\u00f77=139, 0\"),\n    @(\"653\u00f75=130, 3\", \"775\u00f74=193, 3\"),\n    @(\"470\u00f74=117, 2\", \"719\u00f73=239, 2\"),\n    @(\"486\u00f73=162, 0\", \"159\u00f76=26, 3\"),\n    @(\"330\u00f75=66, 0\", \"309\u00f79=34, 3\"),\n    @(\"319\u00f74=79, 3\", \"843\u00f78=105, 3\"),\n    @(\"532\u00f79=59, 1\", \"361\u00f74=90, 1\"),\n    @(\"269\u00f79=29, 8\", \"953\u00f78=119, 1\"),\n    @(\"444\u00f74=111, 0\", \"763\u00f79=84, 7\"),\n    @(\"389\u00f75=77, 4\", \"135\u00f72=67, 1\"),\n    @(\"899\u00f77=128, 3\", \"549\u00f72=274, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
